$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a text value on a cell while preventing Excel from
# auto-coercing numeric-looking strings (e.g. "1.00" -> 1) or losing
# formatting, and without leaving a stray cell style behind.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "67.131.49"
$ws.Range("E2").Value = "  -0.96%  "
Set-TextValue "D3" "2.467.67"
$ws.Range("E3").Value = "  -1.11%  "
Set-TextValue "D4" "1.00"
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "583.41"
$ws.Range("E5").Value = "  -1.46%  "
Set-TextValue "D6" "167.25"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -2.08%  "
Set-TextValue "D9" "2.464.40"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -4.30%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("E14").Value = "  -0.58%  "
Set-TextValue "D16" "66.900.63"
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("E17").Value = "  -4.56%  "
Set-TextValue "D18" "2.431.48"
$ws.Range("E18").Value = "  -2.52%  "
Set-TextValue "D19" "11.37"
$ws.Range("E19").Value = "  -2.96%  "
Set-TextValue "D20" "7.65"
$ws.Range("E20").Value = "  -4.58%  "
Set-TextValue "D21" "355.29"
$ws.Range("E21").Value = "  -2.73%  "
Set-TextValue "D22" "4.02"
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("E23").Value = "  +0.16%  "
Set-TextValue "D24" "69.56"
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("E25").Value = "  -7.20%  "
$ws.Range("E26").Value = "  -7.58%  "
Set-TextValue "D27" "9.00"
$ws.Range("E27").Value = "  -8.79%  "
Set-TextValue "D28" "0.998"
$ws.Range("E28").Value = "  -0.07%  "
Set-TextValue "D29" "2.592.30"
$ws.Range("E29").Value = "  -0.69%  "
Set-TextValue "D30" "0.0₃0900"
$ws.Range("E30").Value = "  -6.78%  "
Set-TextValue "D31" "510.61"
$ws.Range("E31").Value = "  -4.49%  "
$ws.Range("E32").Value = "  -5.39%  "
Set-TextValue "D33" "1.79"
$ws.Range("E34").Value = "  -5.79%  "
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -7.07%  "
Set-TextValue "D37" "158.57"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D38" "18.58"
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D39" "18.44"
$ws.Range("E39").Value = "  -0.89%  "
$ws.Range("E40").Value = "  -6.06%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D41" "1.67"
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D42" "0.326"
$ws.Range("E42").Value = "  -6.50%  "
$ws.Range("E43").Value = "  -6.49%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D44" "2.32"
$ws.Range("E44").Value = "  -7.32%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D45" "38.62"
$ws.Range("E45").Value = "  -2.98%  "
Set-TextValue "D46" "141.45"
$ws.Range("E46").Value = "  -2.32%  "
Set-TextValue "D47" "3.48"
$ws.Range("E47").Value = "  -5.48%  "
Set-TextValue "D48" "0.517"
$ws.Range("E48").Value = "  -5.62%  "
$ws.Range("E49").Value = "  -6.81%  "
$ws.Range("E50").Value = "  -6.04%  "
$ws.Range("E51").Value = "  -1.96%  "
